$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): bold a few day headers (Dag 1, Dag 3, Dag 5) ---
$ws.Range("I2").Font.Bold = $true
$ws.Range("K2").Font.Bold = $true
$ws.Range("M2").Font.Bold = $true

# --- Row 4 values (task: Skapa grundläggande menylayout ...) ---
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("N4").Value = 6

# --- Row 5 values (task: Skapa medlemsklass ...) ---
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 2

# --- Row 6 values (task: Registreringsfunktion ...) ---
$ws.Range("J6").Value = 5

# --- Row 9 ("Uppskattning" totals row): bold label + blank bold cell ---
$ws.Range("B9").Font.Bold = $true
$ws.Range("C9").Font.Bold = $true
$ws.Range("J9").Value = 11
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = 2

# --- Row 10 ("Idealt" totals row): bold label + blank bold cell ---
$ws.Range("B10").Font.Bold = $true
$ws.Range("C10").Font.Bold = $true

# --- Update the active selection to K6 ---
$ws.Range("K6").Select()
